# Update the two-digit multiplication problems in the document.
# Each old "A×B=" expression is replaced with the new one via Find/Replace.

$d = $word.ActiveDocument

$replacements = @(
    @("89×88=", "19×16="),
    @("60×85=", "58×20="),
    @("95×57=", "74×70="),
    @("57×41=", "14×75="),
    @("41×53=", "37×65="),
    @("69×37=", "31×14="),
    @("17×74=", "44×41="),
    @("52×44=", "13×64="),
    @("63×41=", "44×70="),
    @("67×66=", "43×56="),
    @("16×29=", "74×59="),
    @("96×31=", "19×19="),
    @("59×28=", "37×61="),
    @("44×81=", "91×40="),
    @("24×85=", "17×51="),
    @("55×11=", "33×54="),
    @("72×74=", "94×67="),
    @("54×44=", "62×31="),
    @("69×52=", "23×66="),
    @("71×34=", "16×53="),
    @("14×74=", "25×58="),
    @("50×99=", "29×66="),
    @("93×11=", "92×64="),
    @("78×72=", "69×25="),
    @("66×38=", "38×64=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}
